$d = $word.ActiveDocument
$d.Content.Font.Bold = $false
